$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.310.83"
$ws.Range("E2").Value = "  +3.47%  "

# Row 3
$ws.Range("D3").Value = "3.128.47"
$ws.Range("E3").Value = "  +1.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'220.44"
$ws.Range("E5").Value = "  +5.00%  "

# Row 6
$ws.Range("D6").Value = "'622.80"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").Value = "'0.379"
$ws.Range("E7").Value = "  +2.01%  "

# Row 8
$ws.Range("D8").Value = "'0.895"
$ws.Range("E8").Value = "  +9.14%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "3.121.96"
$ws.Range("E10").Value = "  +1.52%  "

# Row 11
$ws.Range("D11").Value = "'0.740"
$ws.Range("E11").Value = "  +25.01%  "

# Row 12
$ws.Range("E12").Value = "  +6.17%  "

# Row 13
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  +6.51%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'33.96"
$ws.Range("E14").Value = "  +7.47%  "

# Row 15
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.42"
$ws.Range("E15").Value = "  +2.54%  "

# Row 16
$ws.Range("D16").Value = "91.112.19"
$ws.Range("E16").Value = "  +3.64%  "

# Row 17
$ws.Range("D17").Value = "3.697.67"
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").Value = "3.140.05"
$ws.Range("E18").Value = "  +2.48%  "

# Row 19
$ws.Range("D19").Value = "'3.80"
$ws.Range("E19").Value = "  +19.03%  "

# Row 20
$ws.Range("D20").Value = "'0.0000229"
$ws.Range("E20").Value = "  +7.73%  "

# Row 21
$ws.Range("D21").Value = "'14.11"
$ws.Range("E21").Value = "  +7.78%  "

# Row 22
$ws.Range("D22").Value = "'432.57"
$ws.Range("E22").Value = "  +2.84%  "

# Row 23
$ws.Range("D23").Value = "'8.81"
$ws.Range("E23").Value = "  +7.73%  "

# Row 24
$ws.Range("D24").Value = "'5.16"
$ws.Range("E24").Value = "  +6.74%  "

# Row 25
$ws.Range("D25").Value = "'6.10"
$ws.Range("E25").Value = "  +12.27%  "

# Row 26
$ws.Range("D26").Value = "'12.34"
$ws.Range("E26").Value = "  +8.67%  "

# Row 27
$ws.Range("D27").Value = "'83.74"
$ws.Range("E27").Value = "  +2.54%  "

# Row 28
$ws.Range("D28").Value = "3.280.81"
$ws.Range("E28").Value = "  +1.44%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").Value = "'0.168"
$ws.Range("E30").Value = "  +9.57%  "

# Row 31
$ws.Range("D31").Value = "'9.02"
$ws.Range("E31").Value = "  +12.40%  "

# Row 32
$ws.Range("D32").Value = "'0.912"
$ws.Range("E32").Value = "  -16.14%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'526.68"
$ws.Range("E33").Value = "  +4.17%  "

# Row 34
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.87"
$ws.Range("E34").Value = "  +8.74%  "

# Row 35
$ws.Range("D35").Value = "'7.30"
$ws.Range("E35").Value = "  +10.56%  "

# Row 36
$ws.Range("D36").Value = "'0.142"
$ws.Range("E36").Value = "  +9.71%  "

# Row 37
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  +6.09%  "

# Row 38
$ws.Range("D38").Value = "'1.86"
$ws.Range("E38").Value = "  +2.69%  "

# Row 39
$ws.Range("D39").Value = "'23.33"
$ws.Range("E39").Value = "  +4.96%  "

# Row 40
$ws.Range("D40").Value = "'22.29"
$ws.Range("E40").Value = "  +0.33%  "

# Row 41
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.378"
$ws.Range("E43").Value = "  +5.40%  "

# Row 44
$ws.Range("D44").Value = "'0.143"
$ws.Range("E44").Value = "  +6.29%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.92"
$ws.Range("E45").Value = "  +5.60%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0765"
$ws.Range("E46").Value = "  +16.17%  "

# Row 47
$ws.Range("D47").Value = "'143.78"
$ws.Range("E47").Value = "  -3.32%  "

# Row 48
$ws.Range("D48").Value = "'44.04"
$ws.Range("E48").Value = "  +1.40%  "

# Row 49
$ws.Range("D49").Value = "'1.30"
$ws.Range("E49").Value = "  +11.54%  "

# Row 50
$ws.Range("E50").Value = "  +25.26%  "

# Row 51
$ws.Range("D51").Value = "'168.02"
$ws.Range("E51").Value = "  +7.36%  "
